$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foaie1")

# --- Update data values (first table: DHT11-1) ---
$ws.Range("B3").Value = 52
$ws.Range("C2").Value = 60
$ws.Range("C3").Value = 45
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 25

# --- Update data values (second table: DHT11-2) ---
$ws.Range("C12").Value = 60
$ws.Range("B13").Value = 54
$ws.Range("C13").Value = 45
$ws.Range("F13").Value = 24
$ws.Range("G13").Value = 25

$excel.CalculateFull()

# --- Update sheet view: clear frozen/scrolled topLeftCell and move selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B12").Select()
